$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.02698669485309508
$ws.Range("C2").Value = 0.4347110584599287
$ws.Range("D2").Value = 0.2582489653893563
$ws.Range("E2").Value = 0.508182019939073
$ws.Range("F2").Value = 0.5213706822786657
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.3926056075262446
$ws.Range("C3").Value = 0.5860682087796297
$ws.Range("D3").Value = 0.6279895078619996
$ws.Range("E3").Value = 0.792457890276827
$ws.Range("F3").Value = 0.7083247097872002
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.658289211196722
$ws.Range("C4").Value = 0.7393275084707053
$ws.Range("D4").Value = 0.8838880653841045
$ws.Range("E4").Value = 0.9401532137817242
$ws.Range("F4").Value = 0.6918831845362214
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.6301452433794297
$ws.Range("C5").Value = 0.6910958835064263
$ws.Range("D5").Value = 0.6770036873424849
$ws.Range("E5").Value = 0.8228023379539493
$ws.Range("F5").Value = 0.5464266070522937
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.5609293215153877
$ws.Range("C6").Value = 0.635969890021205
$ws.Range("D6").Value = 0.5954006161806298
$ws.Range("E6").Value = 0.7716220682307043
$ws.Range("F6").Value = 0.5484643292655379
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.497781246552584
$ws.Range("C7").Value = 0.594007894527335
$ws.Range("D7").Value = 0.4852508987386884
$ws.Range("E7").Value = 0.6965995253649606
$ws.Range("F7").Value = 0.5056987709686329
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.4027147906846605
$ws.Range("C8").Value = 0.5256491137464235
$ws.Range("D8").Value = 0.3760707213903638
$ws.Range("E8").Value = 0.6132460528942391
$ws.Range("F8").Value = 0.4813686861270564
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.3530114966695706
$ws.Range("C9").Value = 0.4763118611471612
$ws.Range("D9").Value = 0.2905691678227247
$ws.Range("E9").Value = 0.5390446807294593
$ws.Range("F9").Value = 0.4320833917452321
$ws.Range("G9").Value = 9

$ws.Range("B10").Value = 0.1773128318588165
$ws.Range("C10").Value = 0.4336731922191779
$ws.Range("D10").Value = 0.2753392244430637
$ws.Range("E10").Value = 0.5247277622187183
$ws.Range("F10").Value = 0.552154172425228
$ws.Range("G10").Value = 5

$ws.Range("B11").Value = 0.002592199968439235
$ws.Range("C11").Value = 0.6934931008693428
$ws.Range("D11").Value = 0.4809394004540528
$ws.Range("E11").Value = 0.6934979455297995
$ws.Range("F11").Value = 0.9807473486615974
